# SkillsDatabase.xlsx - add the new elemental-ability skill rows (Water,
# Plant, Earth, Wind, Thunder, Light and Darkness) to the Plan1 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row layout, in sheet order (row 5 .. row 11).
$rows = @(
    @{ Row = 5;  Id = 3; Name = "DoubleDelta" },
    @{ Row = 6;  Id = 4; Name = "ThornyThrust" },
    @{ Row = 7;  Id = 5; Name = "RockyRumble" },
    @{ Row = 8;  Id = 6; Name = "ChargedContact" },
    @{ Row = 9;  Id = 7; Name = "GulibleGust" },
    @{ Row = 10; Id = 8; Name = "LuminousLash" },
    @{ Row = 11; Id = 9; Name = "ShadowSlice" }
)

# First, stamp the formatting down for every new row by duplicating the
# row right above it (copy + paste-format-only) - this reuses the existing
# "Bom" (green) cell style shared by every data row instead of minting a
# brand-new style entry.
foreach ($r in $rows) {
    $srcRow = $r.Row - 1
    $ws.Range("A" + $srcRow + ":E" + $srcRow).Copy() | Out-Null
    $ws.Range("A" + $r.Row + ":E" + $r.Row).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# Fill in the skill names (column B) in the same order they were typed by
# the author - row 5, then row 7, then row 6, then rows 8-11 - so new
# unique strings land in the shared-string table in that same order.
$ws.Range("B5").Value = "DoubleDelta"
$ws.Range("B7").Value = "RockyRumble"
$ws.Range("B6").Value = "ThornyThrust"
$ws.Range("B8").Value = "ChargedContact"
$ws.Range("B9").Value = "GulibleGust"
$ws.Range("B10").Value = "LuminousLash"
$ws.Range("B11").Value = "ShadowSlice"

# Fill in the remaining columns (ID, Type, Mana, CoolDown) row by row.
foreach ($r in $rows) {
    $ws.Range("A" + $r.Row).Value = $r.Id
    $ws.Range("C" + $r.Row).Value = "DamageSkill"
    $ws.Range("D" + $r.Row).Value = 5
    $ws.Range("E" + $r.Row).Value = 2
}

# Widen column B slightly so the longer skill names keep fitting.
$ws.Columns("B").ColumnWidth = 14.85546875

# Match the author's final selection/view state.
$ws.Range("D13").Select() | Out-Null
